$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66; this shifts existing rows 66.. down by one,
# producing the new row 167 (old row 166 content, unchanged) and making
# room for fresh data at row 66.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new data.
$ws.Cells.Item(66, 1).Value = 8
$ws.Cells.Item(66, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(66, 3).Value = "Coquimbo"
$ws.Cells.Item(66, 4).Value = 45272
$ws.Cells.Item(66, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(66, 5).Value = 4
$ws.Cells.Item(66, 6).Value = "Fruta"
$ws.Cells.Item(66, 7).Value = 100109
$ws.Cells.Item(66, 8).Value = "Uva"
$ws.Cells.Item(66, 9).Value = 100109001
$ws.Cells.Item(66, 10).Value = "Uva"
$ws.Cells.Item(66, 11).Value = "Flame Seedless"
$ws.Cells.Item(66, 12).Value = "Primera"
$ws.Cells.Item(66, 13).Value = 600
$ws.Cells.Item(66, 14).Value = 11500
$ws.Cells.Item(66, 15).Value = 12000
$ws.Cells.Item(66, 16).Value = 11750
$ws.Cells.Item(66, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(66, 18).Value = "Provincia de Copiapó"
$ws.Cells.Item(66, 19).Value = 1175
$ws.Cells.Item(66, 20).Value = 10
